# Reorders/updates the "Estado de Cuenta" table rows (C:F, rows 16-42) so
# that all periods for DERLY ZARATE LLERENA (45517938) come first, in
# descending period order (2101, 2012 .. 2001), followed by all periods
# for GABRIEL JAIME PAREJA (71729664) in the same descending order, ending
# with his 1912 period. Values in column F (Valor Mora) follow the row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 16; Doc = "45517938";  Name = "DERLY ZARATE LLERENA";  Periodo = "2101"; Valor = 26500 },
    @{ Row = 17; Doc = "45517938";  Name = "DERLY ZARATE LLERENA";  Periodo = "2012"; Valor = 33125 },
    @{ Row = 18; Doc = "45517938";  Name = "DERLY ZARATE LLERENA";  Periodo = "2011"; Valor = 33125 },
    @{ Row = 19; Doc = "45517938";  Name = "DERLY ZARATE LLERENA";  Periodo = "2010"; Valor = 33125 },
    @{ Row = 20; Doc = "45517938";  Name = "DERLY ZARATE LLERENA";  Periodo = "2009"; Valor = 33125 },
    @{ Row = 21; Doc = "45517938";  Name = "DERLY ZARATE LLERENA";  Periodo = "2008"; Valor = 33125 },
    @{ Row = 22; Doc = "45517938";  Name = "DERLY ZARATE LLERENA";  Periodo = "2007"; Valor = 33125 },
    @{ Row = 23; Doc = "45517938";  Name = "DERLY ZARATE LLERENA";  Periodo = "2006"; Valor = 33125 },
    @{ Row = 24; Doc = "45517938";  Name = "DERLY ZARATE LLERENA";  Periodo = "2005"; Valor = 33125 },
    @{ Row = 25; Doc = "45517938";  Name = "DERLY ZARATE LLERENA";  Periodo = "2004"; Valor = 33125 },
    @{ Row = 26; Doc = "45517938";  Name = "DERLY ZARATE LLERENA";  Periodo = "2003"; Valor = 33125 },
    @{ Row = 27; Doc = "45517938";  Name = "DERLY ZARATE LLERENA";  Periodo = "2002"; Valor = 33125 },
    @{ Row = 28; Doc = "45517938";  Name = "DERLY ZARATE LLERENA";  Periodo = "2001"; Valor = 33125 },
    @{ Row = 29; Doc = "71729664";  Name = "GABRIEL JAIME PAREJA";  Periodo = "2101"; Valor = 26500 },
    @{ Row = 30; Doc = "71729664";  Name = "GABRIEL JAIME PAREJA";  Periodo = "2012"; Valor = 33125 },
    @{ Row = 31; Doc = "71729664";  Name = "GABRIEL JAIME PAREJA";  Periodo = "2011"; Valor = 33125 },
    @{ Row = 32; Doc = "71729664";  Name = "GABRIEL JAIME PAREJA";  Periodo = "2010"; Valor = 33125 },
    @{ Row = 33; Doc = "71729664";  Name = "GABRIEL JAIME PAREJA";  Periodo = "2009"; Valor = 33125 },
    @{ Row = 34; Doc = "71729664";  Name = "GABRIEL JAIME PAREJA";  Periodo = "2008"; Valor = 33125 },
    @{ Row = 35; Doc = "71729664";  Name = "GABRIEL JAIME PAREJA";  Periodo = "2007"; Valor = 33125 },
    @{ Row = 36; Doc = "71729664";  Name = "GABRIEL JAIME PAREJA";  Periodo = "2006"; Valor = 33125 },
    @{ Row = 37; Doc = "71729664";  Name = "GABRIEL JAIME PAREJA";  Periodo = "2005"; Valor = 33125 },
    @{ Row = 38; Doc = "71729664";  Name = "GABRIEL JAIME PAREJA";  Periodo = "2004"; Valor = 33125 },
    @{ Row = 39; Doc = "71729664";  Name = "GABRIEL JAIME PAREJA";  Periodo = "2003"; Valor = 33125 },
    @{ Row = 40; Doc = "71729664";  Name = "GABRIEL JAIME PAREJA";  Periodo = "2002"; Valor = 33125 },
    @{ Row = 41; Doc = "71729664";  Name = "GABRIEL JAIME PAREJA";  Periodo = "2001"; Valor = 33125 },
    @{ Row = 42; Doc = "71729664";  Name = "GABRIEL JAIME PAREJA";  Periodo = "1912"; Valor = 19875 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.Doc
    $ws.Cells.Item($r.Row, 4).Value = $r.Name
    $ws.Cells.Item($r.Row, 5).Value = $r.Periodo
    $ws.Cells.Item($r.Row, 6).Value = $r.Valor
}
